# Apply the "Mid Paper 1" percentage column + A-level grading fix.
#
# Summary of the change:
#  - A new column E ("Mid Paper 1") is added to both worksheets, holding the
#    percentage-of-total mark (calculated against the grading total).
#  - On "Senior Six" (sheet 1) the existing "Paper 1" marks already used the
#    correct grading total, so column E simply mirrors column D.
#  - On "Senior Five" (sheet 2) the "Paper 1" marks were graded against the
#    wrong total; column E carries the corrected percentage figures (the
#    same corrected figures as "Senior Six"), fixing the grading mismatch
#    described in the commit message.
#  - The previously-active sheet/selection moves from "Senior Six" (A1-ish)
#    to "Senior Five", which becomes the active tab with E7 selected.

$wb = $excel.ActiveWorkbook

$seniorSix  = $wb.Worksheets.Item(1)
$seniorFive = $wb.Worksheets.Item(2)

# --- Senior Six ("Senior Six" sheet): add "Mid Paper 1" header + values ---
$seniorSix.Range("E1").Value = "Mid Paper 1"

$seniorSix.Range("E2").Value = 67
$seniorSix.Range("E3").Value = 48
$seniorSix.Range("E4").Value = 49
$seniorSix.Range("E5").Value = 50
$seniorSix.Range("E6").Value = 51

# --- Senior Five ("Senior Five" sheet): add "Mid Paper 1" header + values ---
$seniorFive.Range("E1").Value = "Mid Paper 1"

$seniorFive.Range("E2").Value = 67
$seniorFive.Range("E3").Value = 48
$seniorFive.Range("E4").Value = 49
$seniorFive.Range("E5").Value = 50
$seniorFive.Range("E6").Value = 51

# --- Restore the cursor on Senior Six to E1, then switch the active tab to
#     Senior Five with E7 selected (matches the saved view state). ---
[void]$seniorSix.Range("E1").Select()

[void]$seniorFive.Activate()
[void]$seniorFive.Range("E7").Select()
